# Update LR-pairs data with new TPM-based values.
# The table now has 20 data rows (rows 2-21) covering all
# Sending cluster x Target cluster combinations (including the
# newly added "Resolving-Mac" target cluster).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Pf4", "Sdc2", "ECs", 1.0, 0.3333333333333333, 0.2859746666666667, 0.857924, 0.001281503322100235, 0.001281503322100236, 3.0, 1.0, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 0.5223516029946667, 4.701164426952, 0.0000358227966424518, 0.0000358227966424518),
    @("ECs", "Pf4", "Sdc2", "FAPs", 1.0, 0.3333333333333333, 0.2859746666666667, 0.857924, 0.001281503322100235, 0.001281503322100236, 3.0, 1.0, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 12.66583504982178, 113.992515448396, 0.0008686211178359045, 0.0008686211178359046),
    @("ECs", "Pf4", "Sdc2", "MuSCs", 1.0, 0.3333333333333333, 0.2859746666666667, 0.857924, 0.001281503322100235, 0.001281503322100236, 3.0, 1.0, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 5.464652537976889, 49.181872841792, 0.0003747650729265647, 0.0003747650729265647),
    @("ECs", "Pf4", "Sdc2", "Resolving-Mac", 1.0, 0.3333333333333333, 0.2859746666666667, 0.857924, 0.001281503322100235, 0.001281503322100236, 2.0, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 0.03345493702977777, 0.301094433268, 0.000002294334695314405, 0.000002294334695314406),
    @("FAPs", "Pf4", "Sdc2", "ECs", 1.0, 0.3333333333333333, 0.8039883333333333, 2.411965, 0.003602814655248594, 0.003602814655248594, 3.0, 1.0, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 1.468537754063333, 13.21683978657, 0.0001007121046895893, 0.0001007121046895894),
    @("FAPs", "Pf4", "Sdc2", "FAPs", 1.0, 0.3333333333333333, 0.8039883333333333, 2.411965, 0.003602814655248594, 0.003602814655248594, 3.0, 1.0, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 35.60869125463722, 320.478221291735, 0.002442038845493398, 0.002442038845493398),
    @("FAPs", "Pf4", "Sdc2", "MuSCs", 1.0, 0.3333333333333333, 0.8039883333333333, 2.411965, 0.003602814655248594, 0.003602814655248594, 3.0, 1.0, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 15.36330800719111, 138.26977206472, 0.001053613419278773, 0.001053613419278773),
    @("FAPs", "Pf4", "Sdc2", "Resolving-Mac", 1.0, 0.3333333333333333, 0.8039883333333333, 2.411965, 0.003602814655248594, 0.003602814655248594, 2.0, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 0.0940551111672222, 0.8464960005049998, 0.000006450285786834276, 0.000006450285786834277),
    @("Inflammatory-Mac", "Pf4", "Sdc2", "ECs", 3.0, 1.0, 54.865331, 164.595993, 0.2458613022061244, 0.2458613022061244, 3.0, 1.0, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 100.215148183346, 901.936333650114, 0.006872740225709293, 0.006872740225709294),
    @("Inflammatory-Mac", "Pf4", "Sdc2", "FAPs", 3.0, 1.0, 54.865331, 164.595993, 0.2458613022061244, 0.2458613022061244, 3.0, 1.0, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 2429.988783621417, 21869.89905259275, 0.1666482758740526, 0.1666482758740526),
    @("Inflammatory-Mac", "Pf4", "Sdc2", "MuSCs", 3.0, 1.0, 54.865331, 164.595993, 0.2458613022061244, 0.2458613022061244, 3.0, 1.0, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 1048.414441009083, 9435.729969081744, 0.07190010924052169, 0.07190010924052169),
    @("Inflammatory-Mac", "Pf4", "Sdc2", "Resolving-Mac", 3.0, 1.0, 54.865331, 164.595993, 0.2458613022061244, 0.2458613022061244, 2.0, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 6.418457323922333, 57.76611591530099, 0.0004401768658408285, 0.0004401768658408285),
    @("MuSCs", "Pf4", "Sdc2", "ECs", 2.0, 0.6666666666666666, 0.04972033333333333, 0.149161, 0.0002228056529806757, 0.0002228056529806757, 3.0, 1.0, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 0.09081747037533333, 0.8173572333779999, 0.00000622824885419309, 0.000006228248854193091),
    @("MuSCs", "Pf4", "Sdc2", "FAPs", 2.0, 0.6666666666666666, 0.04972033333333333, 0.149161, 0.0002228056529806757, 0.0002228056529806757, 3.0, 1.0, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 2.202116529979889, 19.819048769819, 0.000151020830000701, 0.000151020830000701),
    @("MuSCs", "Pf4", "Sdc2", "MuSCs", 2.0, 0.6666666666666666, 0.04972033333333333, 0.149161, 0.0002228056529806757, 0.0002228056529806757, 3.0, 1.0, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 0.9500993528764444, 8.550894175887999, 0.00006515767485558081, 0.00006515767485558081),
    @("MuSCs", "Pf4", "Sdc2", "Resolving-Mac", 2.0, 0.6666666666666666, 0.04972033333333333, 0.149161, 0.0002228056529806757, 0.0002228056529806757, 2.0, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 0.005816566341888888, 0.05234909707699999, 0.0000003988992702008476, 0.0000003988992702008477),
    @("Resolving-Mac", "Pf4", "Sdc2", "ECs", 3.0, 1.0, 167.1506043333333, 501.451813, 0.7490315741635462, 0.749031574163546, 3.0, 1.0, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 305.3116107547194, 2747.804496792474, 0.02093822567393822, 0.02093822567393822),
    @("Resolving-Mac", "Pf4", "Sdc2", "FAPs", 3.0, 1.0, 167.1506043333333, 501.451813, 0.7490315741635462, 0.749031574163546, 3.0, 1.0, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 7403.110239242726, 66627.99215318453, 0.5077042189621703, 0.5077042189621703),
    @("Resolving-Mac", "Pf4", "Sdc2", "MuSCs", 3.0, 1.0, 167.1506043333333, 501.451813, 0.7490315741635462, 0.749031574163546, 3.0, 1.0, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 3194.059057193368, 28746.53151474031, 0.2190481036410021, 0.2190481036410021),
    @("Resolving-Mac", "Pf4", "Sdc2", "Resolving-Mac", 3.0, 1.0, 167.1506043333333, 501.451813, 0.7490315741635462, 0.749031574163546, 2.0, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 19.55422488167122, 175.988023935041, 0.001341025886435408, 0.001341025886435408)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
